# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-26 19:12:06
#
# The upstream edit reorders the two comma-separated entries held in the
# "Recorded By" column (column G) for a specific set of session rows on the
# active worksheet - e.g. "System, dnasr281@gmail.com" becomes
# "dnasr281@gmail.com, System". Only rows that had exactly two
# comma-separated authors in column G were touched by the sync; rows with a
# single author (e.g. just "dnasr281@gmail.com") or three authors
# (e.g. "backup@backdoor.com, system, System") were left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (in the "Recorded By" / column G) whose two comma-separated values
# need to swap places, exactly as produced by the upstream sync.
$rowsToSwap = @(
    3,4,6,7,10,11,12,13,14,15,17,18,19,20,21,22,24,26,29,30,
    32,33,36,37,38,39,40,41,43,44,45,46,47,48,50,52,55,56,58,59,
    62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,87,90,
    92,93,94,96,99,101,109,110,111,112,113,116,118,119,120,122,125,127,
    135,136,137,138,139,142,144,145,146,148,151,153
)

foreach ($r in $rowsToSwap) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7 ("Recorded By")
    $current = $cell.Value()
    if ($null -ne $current) {
        $parts = $current -split ",\s*"
        if ($parts.Count -eq 2) {
            $swapped = $parts[1].Trim() + ", " + $parts[0].Trim()
            $cell.Value = $swapped
        }
    }
}
